$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Лист6")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "Лист7"
